$wb = $excel.ActiveWorkbook

# --- Update the text block in "Hoja1"!A1 with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
[string]$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 5.15 = 20427.84 pesos"), "1000 Bs = 5.06 = 20055.68 pesos"
$text = $text -replace [regex]::Escape("20427.84 pesos = 5.14 = 967.42 Bs"), "20055.68 pesos = 5.02 = 953.36 Bs"
$cell.Value = $text

# --- Update the rate figures on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 197.55
$ws2.Range("O10").Value = 3962
$ws2.Range("N12").Value = 3997
$ws2.Range("O12").Value = 190
